$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)
$shape.TextFrame.TextRange.Text = "What is it ??"
